$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# "strain" sheet: every row (1-8) now has the same strain layout:
#   A:C = MG1655, D:F = UV5, G:I = WT, J:L = 3_19 (date-ish number format carried over)
# ------------------------------------------------------------------
$wsStrain = $wb.Worksheets.Item("strain")

for ($r = 1; $r -le 8; $r++) {
    $wsStrain.Range("A${r}:C${r}").Value = "MG1655"
    $wsStrain.Range("D${r}:F${r}").Value = "UV5"
    $wsStrain.Range("G${r}:I${r}").Value = "WT"
    $wsStrain.Range("J${r}:L${r}").Value = "3_19"
}

# D2:F8 used to carry the "d-mmm" number format (leftover from the old "3_19" column);
# that formatting now belongs to J:L instead.
$wsStrain.Range("D2:F8").Style = "Normal"
$wsStrain.Range("J1:L8").NumberFormat = "d-mmm"

# ------------------------------------------------------------------
# "pos_selection" sheet: tc concentration gradient changes from
# 0 / 2.13 / 4.57 / 10 / 21.3 / 45.7 / 100  ->  0 / 10 / 25 / 40 / 55 / 70 / 85 / 100
# ------------------------------------------------------------------
$wsPos = $wb.Worksheets.Item("pos_selection")

$wsPos.Range("A1:L1").Value = "0_µg/ml_tc"
$wsPos.Range("A2:L2").Value = "10_µg/ml_tc"
$wsPos.Range("A3:L3").Value = "25_µg/ml_tc"
$wsPos.Range("A4:L4").Value = "40_µg/ml_tc"
$wsPos.Range("A5:L5").Value = "55_µg/ml_tc"
$wsPos.Range("A6:L6").Value = "70_µg/ml_tc"
$wsPos.Range("A7:L7").Value = "85_µg/ml_tc"
$wsPos.Range("A8:L8").Value = "100_µg/ml_tc"

# ------------------------------------------------------------------
# Selection / active-tab bookkeeping: "pos_selection" loses its
# tabSelected flag and its selection collapses to the full used range,
# "strain" becomes the active sheet (tabSelected + the live cell
# selection at L17).
# ------------------------------------------------------------------
$wsPos.Range("A1:L8").Select()
$wsStrain.Activate()
$wsStrain.Range("L17").Select()
